# "Copycat completed and Suppression attack started"
# Populate a new column F with a second batch of random seed values
# (mirrors the existing A/D columns) and highlight it with a new
# light-green fill color, same as selecting F1:F10, typing the values,
# then applying a custom Fill Color (RGB 194, 224, 174) from the
# "More Colors..." picker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$seedValues = @(765765, 293847, 876976, 124897, 111685, 549849, 686986, 979845, 916619, 398097)

for ($i = 0; $i -lt $seedValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 6).Value = $seedValues[$i]
}

# RGB(194, 224, 174) -> 0xAEE0C2 in Excel's BGR-packed OLE_COLOR form
$ws.Range("F1:F10").Interior.Color = 11460802

# Leave the cursor where the author's session ended up
$ws.Range("E27").Select() | Out-Null
